# Fruta / hortaliza, semanal
# Reorders the weekly price records (rows 2-12) so that each row's
# Fecha/Volumen/Precio mínimo/Precio máximo/Precio promedio ponderado/Precio $/Kg
# values match the data of a different original row (the data itself is
# unchanged, only redistributed across the date rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44455
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 950
$ws.Range("P2").Value = 950

$ws.Range("D3").Value = 44442
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 850
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 875
$ws.Range("P3").Value = 875

$ws.Range("D4").Value = 44341
$ws.Range("J4").Value = 1300
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("P4").Value = 950

$ws.Range("D5").Value = 44175
$ws.Range("J5").Value = 1600
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1100
$ws.Range("P5").Value = 1100

$ws.Range("D6").Value = 44284
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 800
$ws.Range("L6").Value = 850
$ws.Range("M6").Value = 825
$ws.Range("P6").Value = 825

$ws.Range("D7").Value = 44449
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 950
$ws.Range("M7").Value = 925
$ws.Range("P7").Value = 925

$ws.Range("D8").Value = 44291
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1100
$ws.Range("P8").Value = 1100

$ws.Range("D9").Value = 44229
$ws.Range("J9").Value = 1500
$ws.Range("K9").Value = 1400
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = 1450
$ws.Range("P9").Value = 1450

$ws.Range("D10").Value = 44407
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 1200
$ws.Range("L10").Value = 1300
$ws.Range("M10").Value = 1250
$ws.Range("P10").Value = 1250

$ws.Range("D11").Value = 44243
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1250
$ws.Range("P11").Value = 1250

$ws.Range("D12").Value = 44453
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 800
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = 850
$ws.Range("P12").Value = 850
